$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full data set (rows 2-10), replacing the previous rows 2-5 (which only
# covered a partial FAPs/MuSCs x Clec11a/Itga11 matrix) with the complete
# ECs/FAPs/MuSCs x ECs/FAPs/MuSCs matrix using updated TPM-derived values.
$data = @(
    @("ECs", "Clec11a", "Itga11", "ECs", 1, 0.3333333333333333, 0.022591, 0.067773, 0.001469689085715816, 0.001469689085715816, 1, 0.3333333333333333, 0.006053333333333333, 0.01816, 0.0002373030766683641, 0.0002373030766683641, 0.0001367508533333333, 0.00123075768, [double]"3.487617417862782e-07", [double]"3.487617417862783e-07"),
    @("ECs", "Clec11a", "Itga11", "FAPs", 1, 0.3333333333333333, 0.022591, 0.067773, 0.001469689085715816, 0.001469689085715816, 3, 1, 24.759128, 74.277384, 0.970608576546119, 0.9706085765461191, 0.559333460648, 5.034001145832, 0.001426492831451995, 0.001426492831451995),
    @("ECs", "Clec11a", "Itga11", "MuSCs", 1, 0.3333333333333333, 0.022591, 0.067773, 0.001469689085715816, 0.001469689085715816, 3, 1, 0.7436886666666668, 2.231066, 0.02915412037721258, 0.02915412037721259, 0.01680067066866667, 0.151206036018, [double]"4.28474925220344e-05", [double]"4.284749252203441e-05"),
    @("FAPs", "Clec11a", "Itga11", "ECs", 3, 1, 14.768619, 44.305857, 0.9607931545923257, 0.9607931545923256, 1, 0.3333333333333333, 0.006053333333333333, 0.01816, 0.0002373030766683641, 0.0002373030766683641, 0.08939937368, 0.80459436312, 0.0002279991716266621, 0.0002279991716266621),
    @("FAPs", "Clec11a", "Itga11", "FAPs", 3, 1, 14.768619, 44.305857, 0.9607931545923257, 0.9607931545923256, 3, 1, 24.759128, 74.277384, 0.970608576546119, 0.9706085765461191, 365.658128204232, 3290.923153838088, 0.9325540761341126, 0.9325540761341126),
    @("FAPs", "Clec11a", "Itga11", "MuSCs", 3, 1, 14.768619, 44.305857, 0.9607931545923257, 0.9607931545923256, 3, 1, 0.7436886666666668, 2.231066, 0.02915412037721258, 0.02915412037721259, 10.983254572618, 98.84929115356202, 0.02801107928658648, 0.02801107928658648),
    @("MuSCs", "Clec11a", "Itga11", "ECs", 3, 1, 0.5800683333333333, 1.740205, 0.03773715632195847, 0.03773715632195847, 1, 0.3333333333333333, 0.006053333333333333, 0.01816, 0.0002373030766683641, 0.0002373030766683641, 0.003511346977777777, 0.0316021228, [double]"8.955143299915752e-06", [double]"8.955143299915753e-06"),
    @("MuSCs", "Clec11a", "Itga11", "FAPs", 3, 1, 0.5800683333333333, 1.740205, 0.03773715632195847, 0.03773715632195847, 3, 1, 24.759128, 74.277384, 0.970608576546119, 0.9706085765461191, 14.36198611374667, 129.25787502372, 0.03662800758055448, 0.03662800758055449),
    @("MuSCs", "Clec11a", "Itga11", "MuSCs", 3, 1, 0.5800683333333333, 1.740205, 0.03773715632195847, 0.03773715632195847, 3, 1, 0.7436886666666668, 2.231066, 0.02915412037721258, 0.02915412037721259, 0.4313902453922223, 3.882512208530001, 0.001100193598104066, 0.001100193598104066)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowValues[$c]
    }
}

$wb.Save()
